# Apply the flume-analysis "pad" column backfill described in the commit
# "updating analysis with new flume load data."
#
# The workbook tracks rainfall-simulator readings; column G ("pad") records
# which flume pad a given measurement belongs to. Several contiguous blocks
# of rows were missing this value. They are backfilled here with the pad id
# that already appears at the boundaries of each block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Backfill column G ("pad") for rows that were missing it ---------------

# Rows 77-109 belong to pad 2 (matches existing G76/G110 markers)
$ws.Range("G77:G109").Value = 2

# Rows 757-777 belong to pad 1 (matches existing G756 marker)
$ws.Range("G757:G777").Value = 1

# Rows 779-856 belong to pad 2 (matches existing G778 marker)
$ws.Range("G779:G856").Value = 2

# Rows 858-948 belong to pad 3 (matches existing G857/G949 markers)
$ws.Range("G858:G948").Value = 3

# --- Fix an isolated missing load value in column F -------------------------
# Row 900 was missing its F (fcontrol) reading; neighboring rows 899 and 901
# both record 3029, so fill the same value in for row 900.
$ws.Range("F900").Value = 3029

# --- Restore the scroll / selection / zoom state recorded in the saved view
$win = $excel.ActiveWindow
$win.ScrollRow = 744
$win.ScrollColumn = 1
$win.Zoom = 115
$ws.Range("G757:G777").Select()
